$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and title date
$ws.Name = "Through 2022-10-26"

# Update the "October (through 10-25)" label to "October (through 10-26)"
$ws.Range("A11").Value = "October (through 10-26)"

# Update October row (row 11) values for years 2015-2022 (columns B-I)
$ws.Range("B11").Value = 25
$ws.Range("C11").Value = 44
$ws.Range("D11").Value = 59
$ws.Range("E11").Value = 56
$ws.Range("F11").Value = 47
$ws.Range("G11").Value = 126
$ws.Range("H11").Value = 162
$ws.Range("I11").Value = 95

# Update Total row (row 12) values for years 2015-2022 (columns B-I)
$ws.Range("B12").Value = 251
$ws.Range("C12").Value = 473
$ws.Range("D12").Value = 686
$ws.Range("E12").Value = 604
$ws.Range("F12").Value = 469
$ws.Range("G12").Value = 1027
$ws.Range("H12").Value = 1409
$ws.Range("I12").Value = 1372
